$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 120
$ws1.Range("F4").Value = 1823
$ws1.Range("F6").Value = 38
$ws1.Range("F7").Value = 82
$ws1.Range("F9").Value = 10321
$ws1.Range("F14").Value = 395
$ws1.Range("F15").Value = 7179
$ws1.Range("F18").Value = 113
$ws1.Range("F20").Value = 257

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 15

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 120
$ws4.Range("F4").Value = 1823
$ws4.Range("F6").Value = 38
$ws4.Range("F7").Value = 15
$ws4.Range("F8").Value = 82
$ws4.Range("F12").Value = 10321
$ws4.Range("F17").Value = 395
$ws4.Range("F18").Value = 7179
$ws4.Range("F21").Value = 113
$ws4.Range("F23").Value = 257
